$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.112.16"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.831.69"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'243.39"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'0.6274"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.07504"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'23.31"
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("D11").Value = "'0.07702"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "1.816.68"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "'5.030"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "'0.6694"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "'0.000009395"
$ws.Range("E16").Value = "  -6.99%  "
$ws.Range("D17").Value = "'5.997"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "29.088.61"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "2.050.96"
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "'223.63"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'7.135"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'160.19"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "'8.517"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'1.493"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'0.05826"
$ws.Range("E30").Value = "  +10.90%  "
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "'4.127"
$ws.Range("E32").Value = "  +2.66%  "
$ws.Range("D33").Value = "'1.211"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("D34").Value = "'0.7419"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "1.235.70"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "'2.763"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").Value = "'6.501"
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("D42").Value = "'0.8942"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'102.23"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000126"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'65.94"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.959.22"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "'0.5090"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "'0.4074"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "'0.07510"
$ws.Range("E50").Value = "  +12.82%  "
$ws.Range("D51").Value = "'8.983"
$ws.Range("E51").Value = "  +1.22%  "
